$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "joint bilateral filter" speed-up ratio column (G = opencv / opencp time).
# G2 is entered first as a standalone formula, then G3:G31 is filled as a
# second operation so Excel groups G3:G31 into one shared-formula block
# (matching how the column was authored interactively).
$ws.Range("G2").Formula = "=B2/C2"
$ws.Range("G3:G31").Formula = "=B3/C3"

# Move/resize the chart to its new anchor position (from col13,row5 to col25,row35).
$co = $ws.ChartObjects(1)
$co.Top = 74.25
$co.Left = 764.1875787401575
$co.Width = 717.75
$co.Height = 408.75

# Update the visible selection to the newly-added column.
$ws.Range("G2:G31").Select()
